$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark from the "2.2 Das Bestellungsarray..." paragraph.
try {
    $old = $d.Bookmarks.Item("_GoBack")
    $old.Delete()
} catch {
    # no pre-existing bookmark - nothing to remove
}

# 2) Locate the "Die Bestellung bestellen lassen" paragraph and insert a new list
#    paragraph right after it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Die Bestellung bestellen lassen") {
        $target = $p
    }
}

$insertPos = $target.Range.End
$insertRange = $d.Range($insertPos, $insertPos)

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
'<w:r><w:t>3.1 Bestellknopf: Session l' + [char]0xF6 + 'schen und Liste erstellen</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> - ERLEDIGT</w:t></w:r>' +
'</w:p>' +
'<w:p/>' +
'</w:body></w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'

$insertRange.InsertXML($newParaXml)

# InsertXML always terminates the inserted fragment with its own paragraph mark,
# so an extra blank paragraph is left over after the new one - remove it.
$newPara = $target.Next()
$trailing = $newPara.Next()
$trailing.Range.Delete()

# 3) Re-add the "_GoBack" bookmark at the end of the freshly inserted paragraph.
$bmRange = $newPara.Range.Duplicate
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)
